$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 76: GenomeWeb article about Chinese MDx / GenePlus IPO
$ws.Range("A76").Value = "https://www.genomeweb.com/cancer/chinese-mdx-precision-medicine-firm-geneplus-files-ipo-hong-kong-stock-exchange"
$ws.Range("B76").Value = "CDx"
$ws.Range("C76").Value = "Chinese MDx, Precision Medicine Firm GenePlus Files for IPO on Hong Kong Stock Exchange"

# New row 77: 360Dx article, same topic/title/keywords
$ws.Range("A77").Value = "https://www.360dx.com/cancer/chinese-mdx-precision-medicine-firm-geneplus-files-ipo-hong-kong-stock-exchange"
$ws.Range("B77").Value = "CDx"
$ws.Range("C77").Value = "Chinese MDx, Precision Medicine Firm GenePlus Files for IPO on Hong Kong Stock Exchange"

# Wire up the hyperlinks on column A exactly like the existing rows
$ws.Hyperlinks.Add($ws.Range("A76"), "https://www.genomeweb.com/cancer/chinese-mdx-precision-medicine-firm-geneplus-files-ipo-hong-kong-stock-exchange")
$ws.Hyperlinks.Add($ws.Range("A77"), "https://www.360dx.com/cancer/chinese-mdx-precision-medicine-firm-geneplus-files-ipo-hong-kong-stock-exchange")

# Match the existing hyperlink cell style used for the rest of column A
$ws.Range("A76").Style = "Hyperlink"
$ws.Range("A77").Style = "Hyperlink"
